$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.079.02"
$ws.Range("E2").Value = "  -0.27%  "

# Row 3
$ws.Range("D3").Value = "2.409.88"
$ws.Range("E3").Value = "  -0.93%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'561.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "

# Row 6
$ws.Range("D6").Value = "'142.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.94%  "

# Row 7
$ws.Range("E7").Value = "  +0.13%  "

# Row 8
$ws.Range("D8").Value = "'0.528"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.57%  "

# Row 9
$ws.Range("E9").Value = "  +0.09%  "

# Row 10
$ws.Range("E10").Value = "  -2.03%  "

# Row 11
$ws.Range("D11").Value = "'5.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.95%  "

# Row 12
$ws.Range("D12").Value = "'0.349"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.58%  "

# Row 13
$ws.Range("D13").Value = "'25.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.45%  "

# Row 14
$ws.Range("D14").Value = "'0.0000173"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.88%  "

# Row 15
$ws.Range("D15").Value = "2.845.25"
$ws.Range("E15").Value = "  +0.21%  "

# Row 16
$ws.Range("D16").Value = "62.105.76"
$ws.Range("E16").Value = "  +0.12%  "

# Row 17
$ws.Range("D17").Value = "2.409.99"
$ws.Range("E17").Value = "  -0.59%  "

# Row 18
$ws.Range("D18").Value = "'11.27"
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'6.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.84%  "

# Row 20
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.37%  "

# Row 21
$ws.Range("D21").Value = "'320.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.39%  "

# Row 22
$ws.Range("E22").Value = "  -0.12%  "

# Row 23
$ws.Range("D23").Value = "'65.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.58%  "

# Row 24
$ws.Range("E24").Value = "  -1.60%  "

# Row 25
$ws.Range("D25").Value = "'8.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.30%  "

# Row 26
$ws.Range("D26").Value = "'568.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.44%  "

# Row 27
$ws.Range("E27").Value = "  +0.33%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0936"
$ws.Range("E29").Value = "  -0.22%  "

# Row 30
$ws.Range("D30").Value = "'8.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.61%  "

# Row 31
$ws.Range("E31").Value = "  -3.32%  "

# Row 32
$ws.Range("E32").Value = "  -0.75%  "

# Row 33
$ws.Range("E33").Value = "  -0.13%  "

# Row 34
$ws.Range("E34").Value = "  -2.67%  "

# Row 35
$ws.Range("E35").Value = "  +0.17%  "

# Row 36
$ws.Range("D36").Value = "'4.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.07%  "

# Row 37
$ws.Range("D37").Value = "'5.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.92%  "

# Row 38
$ws.Range("E38").Value = "  -1.17%  "

# Row 39
$ws.Range("D39").Value = "'151.65"
$ws.Range("D39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'18.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.99%  "

# Row 41
$ws.Range("E41").Value = "  -9.94%  "

# Row 42
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "

# Row 43
$ws.Range("E43").Value = "  -0.15%  "

# Row 44
$ws.Range("D44").Value = "'147.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.62%  "

# Row 45
$ws.Range("D45").Value = "'3.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.79%  "

# Row 46
$ws.Range("D46").Value = "'0.0531"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.18%  "

# Row 47
$ws.Range("D47").Value = "'19.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.84%  "

# Row 48
$ws.Range("E48").Value = "  -0.22%  "

# Row 49
$ws.Range("E49").Value = "  +0.40%  "

# Row 50
$ws.Range("E50").Value = "  -1.33%  "

# Row 51
$ws.Range("E51").Value = "  +0.15%  "
